$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2210", "_new" -> "_FV2304" -------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the header+data range into a real Excel Table ---------------
# Stash the header row's existing look (bold/shaded/bordered) in a scratch
# row, strip formatting before the table is created (so Excel doesn't bake
# a header-row differential format / dxf into the table definition), then
# restore the original look afterwards.
$headerRange = $ws.Range("A1:U1")
$backupRange = $ws.Range("A100:U100")

$headerRange.Copy()
$backupRange.PasteSpecial(-4122)  # xlPasteFormats
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"

$backupRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$backupRange.Clear()
$excel.CutCopyMode = $false

$tbl.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
